# 12-17-2024 job postings update
# The "dayofPosting" bucket (column F) advances by one step for the
# last row(s) of several buckets, i.e. jobs that have aged into the
# next reporting bucket since the previous export:
#   2D -> 3D (row 5)
#   3D -> 4D (row 11)
#   4D -> 5D (rows 17-18)
#   5D -> 6D (row 25)
#   6D -> 1W (rows 36-37)
#   1W -> 2W (rows 89-90)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "3D"
$ws.Range("F11").Value = "4D"
$ws.Range("F17").Value = "5D"
$ws.Range("F18").Value = "5D"
$ws.Range("F25").Value = "6D"
$ws.Range("F36").Value = "1W"
$ws.Range("F37").Value = "1W"
$ws.Range("F89").Value = "2W"
$ws.Range("F90").Value = "2W"
